# Updated symbol list data: Price (D), Volume 1h (E), and Hora (G) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '308.64'
Set-TextValue $ws.Range("E2") '0.30%'
Set-TextValue $ws.Range("G2") '9'

Set-TextValue $ws.Range("D3") '40.75'
Set-TextValue $ws.Range("E3") '1.76%'
Set-TextValue $ws.Range("G3") '9'

Set-TextValue $ws.Range("D4") '5.119'
Set-TextValue $ws.Range("E4") '-0.31%'
Set-TextValue $ws.Range("G4") '9'

Set-TextValue $ws.Range("D5") '0.07620'
Set-TextValue $ws.Range("E5") '-1.63%'
Set-TextValue $ws.Range("G5") '9'

Set-TextValue $ws.Range("D6") '1.608'
Set-TextValue $ws.Range("E6") '-1.18%'
Set-TextValue $ws.Range("G6") '9'

Set-TextValue $ws.Range("D7") '0.9041'
Set-TextValue $ws.Range("E7") '2.63%'
Set-TextValue $ws.Range("G7") '9'

Set-TextValue $ws.Range("E8") '0.27%'
Set-TextValue $ws.Range("G8") '9'

Set-TextValue $ws.Range("D9") '0.1110'
Set-TextValue $ws.Range("E9") '9.64%'
Set-TextValue $ws.Range("G9") '9'

Set-TextValue $ws.Range("D10") '0.1783'
Set-TextValue $ws.Range("E10") '1.95%'
Set-TextValue $ws.Range("G10") '9'

Set-TextValue $ws.Range("E11") '2.30%'
Set-TextValue $ws.Range("G11") '9'

Set-TextValue $ws.Range("D12") '0.04155'
Set-TextValue $ws.Range("E12") '-5.54%'
Set-TextValue $ws.Range("G12") '9'

Set-TextValue $ws.Range("E13") '-0.21%'
Set-TextValue $ws.Range("G13") '9'

Set-TextValue $ws.Range("D14") '0.001260'
Set-TextValue $ws.Range("E14") '0.17%'
Set-TextValue $ws.Range("G14") '9'

Set-TextValue $ws.Range("D15") '0.005870'
Set-TextValue $ws.Range("E15") '0.57%'
Set-TextValue $ws.Range("G15") '9'

Set-TextValue $ws.Range("D16") '3.354'
Set-TextValue $ws.Range("G16") '9'

Set-TextValue $ws.Range("D17") '4.250'
Set-TextValue $ws.Range("E17") '0.07%'
Set-TextValue $ws.Range("G17") '9'

Set-TextValue $ws.Range("E18") '-0.95%'
Set-TextValue $ws.Range("G18") '9'

Set-TextValue $ws.Range("D19") '6.535'
Set-TextValue $ws.Range("E19") '-6.92%'
Set-TextValue $ws.Range("G19") '9'

Set-TextValue $ws.Range("D20") '0.1365'
Set-TextValue $ws.Range("E20") '1.93%'
Set-TextValue $ws.Range("G20") '9'

Set-TextValue $ws.Range("E21") '1.91%'
Set-TextValue $ws.Range("G21") '9'

Set-TextValue $ws.Range("D22") '0.04065'
Set-TextValue $ws.Range("E22") '-2.35%'
Set-TextValue $ws.Range("G22") '9'

Set-TextValue $ws.Range("D23") '0.001230'
Set-TextValue $ws.Range("E23") '2.40%'
Set-TextValue $ws.Range("G23") '9'

Set-TextValue $ws.Range("D24") '0.004121'
Set-TextValue $ws.Range("E24") '0.41%'
Set-TextValue $ws.Range("G24") '9'

Set-TextValue $ws.Range("D25") '0.0001300'
Set-TextValue $ws.Range("E25") '-0.06%'
Set-TextValue $ws.Range("G25") '9'

Set-TextValue $ws.Range("G26") '9'

Set-TextValue $ws.Range("G27") '9'

Set-TextValue $ws.Range("G28") '9'

Set-TextValue $ws.Range("G29") '9'

Set-TextValue $ws.Range("G30") '9'

Set-TextValue $ws.Range("G31") '9'

Set-TextValue $ws.Range("G32") '9'

Set-TextValue $ws.Range("G33") '9'

Set-TextValue $ws.Range("G34") '9'

Set-TextValue $ws.Range("G35") '9'

Set-TextValue $ws.Range("G36") '9'

Set-TextValue $ws.Range("G37") '9'

Set-TextValue $ws.Range("D38") '0.02420'
Set-TextValue $ws.Range("E38") '2.47%'
Set-TextValue $ws.Range("G38") '9'

Set-TextValue $ws.Range("D39") '0.05191'
Set-TextValue $ws.Range("E39") '-0.12%'
Set-TextValue $ws.Range("G39") '9'

Set-TextValue $ws.Range("D40") '0.007744'
Set-TextValue $ws.Range("E40") '-2.56%'
Set-TextValue $ws.Range("G40") '9'

Set-TextValue $ws.Range("D41") '0.1304'
Set-TextValue $ws.Range("E41") '-1.89%'
Set-TextValue $ws.Range("G41") '9'

Set-TextValue $ws.Range("D42") '0.007008'
Set-TextValue $ws.Range("E42") '10.11%'
Set-TextValue $ws.Range("G42") '9'

Set-TextValue $ws.Range("D43") '0.001950'
Set-TextValue $ws.Range("E43") '-1.14%'
Set-TextValue $ws.Range("G43") '9'

Set-TextValue $ws.Range("D44") '0.008782'
Set-TextValue $ws.Range("E44") '-1.68%'
Set-TextValue $ws.Range("G44") '9'

Set-TextValue $ws.Range("D45") '0.3333'
Set-TextValue $ws.Range("E45") '0.20%'
Set-TextValue $ws.Range("G45") '9'

Set-TextValue $ws.Range("D46") '0.00006932'
Set-TextValue $ws.Range("E46") '5.49%'
Set-TextValue $ws.Range("G46") '9'

Set-TextValue $ws.Range("E47") '-0.06%'
Set-TextValue $ws.Range("G47") '9'

Set-TextValue $ws.Range("D48") '0.03118'
Set-TextValue $ws.Range("E48") '389.76%'
Set-TextValue $ws.Range("G48") '9'

Set-TextValue $ws.Range("D49") '0.004200'
Set-TextValue $ws.Range("E49") '-40.05%'
Set-TextValue $ws.Range("G49") '9'

Set-TextValue $ws.Range("D50") '0.00002100'
Set-TextValue $ws.Range("E50") '-0.06%'
Set-TextValue $ws.Range("G50") '9'

Set-TextValue $ws.Range("D51") '0.0002000'
Set-TextValue $ws.Range("E51") '-0.06%'
Set-TextValue $ws.Range("G51") '9'
